$wb = $excel.ActiveWorkbook

# A new handback has been generated, updating the "Correspond Handoff Datetime"
# and "Correspond Handback DateTime" for the 8d52b95f... file in both the
# zh-cn and de-de language sheets. The Overview sheet's
# "Latest HO Xliff Generate Date" column recalculates automatically from
# these values.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-07 05:50:40"
$wsZhCn.Range("K2").Value = "2016-09-07 05:51:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-07 05:50:51"
$wsDeDe.Range("K2").Value = "2016-09-07 05:51:54"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 05:50:51"
